$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.331.95"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.931.42"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.56"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7149"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3260"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.41"
$ws.Range("E9").Value = "  +3.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07186"
$ws.Range("E10").Value = "  +5.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7980"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08089"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("D13").Value = "1.926.52"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.422"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.72"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.83"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "30.308.72"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "251.35"
$ws.Range("E18").Value = "  -3.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008140"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").Value = "2.182.25"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.909"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.717"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.78"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.21"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.312"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1280"
$ws.Range("E29").Value = "  -4.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.360"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.546"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.427"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.197"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05204"
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.266"
$ws.Range("E35").Value = "  +5.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7468"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.760"
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01954"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.89"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.416"
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4522"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.024"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8398"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.75"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.760"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.400"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.58"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06099"
$ws.Range("E50").Value = "  +3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4170"
$ws.Range("E51").Value = "  +1.51%  "
